$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh (GitHub Actions scrape update)
# Column D values are plain text (e.g. "35.925.97" is not a valid number),
# so they are written with a leading apostrophe to force text entry and
# preserve exact formatting (trailing zeros, multi-dot separators, etc).
$ws.Range("D2").Value = "'35.925.97"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "'1.894.07"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'247.23"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "'0.692"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'43.43"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").Value = "'57.13"
$ws.Range("E9").Value = "  +8.56%  "
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("E11").Value = "  +2.80%  "
$ws.Range("D12").Value = "'0.0986"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "'14.95"
$ws.Range("E13").Value = "  +14.43%  "
$ws.Range("E14").Value = "  +6.60%  "
$ws.Range("D15").Value = "'2.169.35"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "'1.899.28"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").Value = "'35.954.78"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "'73.33"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").Value = "'247.00"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "'13.05"
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("E23").Value = "  +5.31%  "
$ws.Range("E24").Value = "  +6.83%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").Value = "'167.63"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("E28").Value = "  +2.67%  "
$ws.Range("D29").Value = "'18.48"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("E31").Value = "  +4.83%  "
$ws.Range("E32").Value = "  +5.38%  "
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "'1.89"
$ws.Range("E34").Value = "  +6.91%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -14.57%  "
$ws.Range("D37").Value = "'0.861"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").Value = "'0.0780"
$ws.Range("E38").Value = "  +14.12%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  +4.72%  "
$ws.Range("D41").Value = "'99.92"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").Value = "'15.31"
$ws.Range("E42").Value = "  +27.23%  "
$ws.Range("D43").Value = "'17.01"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("D45").Value = "'1.318.18"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").Value = "'2.34"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'0.0811"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "'2.76"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "'6.33"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").Value = "'43.05"
$ws.Range("E51").Value = "  -0.44%  "
